$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update part descriptions: change USB mini to micro (row 6 / ref J1)
$ws.Range("C6").Value = "USB_B_Micro"
$ws.Range("D6").Value = "USB_B_Micro"
$ws.Range("F6").Value = "USB Micro Type B connector"

# Widen the Footprint column (E) to fit the updated text
$ws.Columns.Item(5).ColumnWidth = 49.5

# Leave the selection on the edited footprint cell
$ws.Range("E6").Select()
